$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.880.13'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '2.906.31'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = "'569.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.29%  '
$ws.Range("D6").Value = "'144.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = "'0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '2.908.12'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").Value = "'6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.52%  '
$ws.Range("E11").Value = '  -1.57%  '
$ws.Range("D12").Value = "'0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = "'32.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '3.387.79'
$ws.Range("E16").Value = '  -1.75%  '
$ws.Range("D17").Value = '61.859.79'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = "'6.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").Value = '2.892.54'
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").Value = "'437.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").Value = "'13.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("D22").Value = "'0.658"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("D23").Value = "'6.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("D24").Value = "'79.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").Value = "'11.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = "'10.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.30%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = "'2.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("E29").Value = '  +9.23%  '
$ws.Range("D30").Value = "'7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").Value = "'2.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").Value = "'2.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.30%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'25.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("D36").Value = "'0.960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.53%  '
$ws.Range("D37").Value = "'5.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("D38").Value = "'49.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").Value = "'2.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.98%  '
$ws.Range("D40").Value = "'1.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").Value = "'8.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("D43").Value = "'40.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = "'0.271"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.75%  '
$ws.Range("D45").Value = '2.694.91'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").Value = "'133.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").Value = "'0.0335"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = "'340.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.79%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = "'21.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.82%  '
